# Update countries & provincias Spain
# - Swap display order (shared-string text) for six country-name pairs
#   whose ranking changed between snapshots.
# - Update the "Datos actualizados" timestamp text.
# - Update the numeric statistics (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Casos criticos, Muertes hoy, Muertes) for the
#   rows whose country data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Swap the country names that changed rank order (column A labels).
#    Each pair below occupies two consecutive data rows; the label text
#    is swapped while the row (rank) position stays the same.
# ---------------------------------------------------------------------
$ws.Range("A39").Value = "Ucrania"
$ws.Range("A40").Value = "Polonia"

$ws.Range("A58").Value = "Armenia"
$ws.Range("A59").Value = "Argelia"

$ws.Range("A143").Value = "Taiwan"
$ws.Range("A144").Value = "Togo"

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Belice"

$ws.Range("A210").Value = "Montserrat"
$ws.Range("A211").Value = "Seychelles"

$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"

# ---------------------------------------------------------------------
# 2) Update the "last updated" timestamp text in A1.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 09:05"

# ---------------------------------------------------------------------
# 3) Update the numeric statistics that changed.
#    Columns: B=Casos totales, C=Nuevos casos, D=Casos activos,
#             E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes
# ---------------------------------------------------------------------
$ws.Range("D12").Value = 165900
$ws.Range("E12").Value = 8989

$ws.Range("B39").Value = 24012
$ws.Range("C39").Value = 340
$ws.Range("D39").Value = 9690
$ws.Range("E39").Value = 13604
$ws.Range("G39").Value = 10
$ws.Range("H39").Value = 718

$ws.Range("B40").Value = 23786
$ws.Range("D40").Value = 11449
$ws.Range("E40").Value = 11273
$ws.Range("H40").Value = 1064

$ws.Range("B58").Value = 9402
$ws.Range("C58").Value = 120
$ws.Range("D58").Value = 3402
$ws.Range("E58").Value = 5861
$ws.Range("G58").Value = 8
$ws.Range("H58").Value = 139

$ws.Range("B59").Value = 9394
$ws.Range("D59").Value = 5748
$ws.Range("E59").Value = 2993
$ws.Range("H59").Value = 653

$ws.Range("B76").Value = 3892
$ws.Range("C76").Value = 16
$ws.Range("D76").Value = 2156
$ws.Range("E76").Value = 1209
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 527

$ws.Range("B126").Value = 794
$ws.Range("C126").Value = 11
$ws.Range("D126").Value = 624
$ws.Range("E126").Value = 158

$ws.Range("B142").Value = 449
$ws.Range("C142").Value = 1
$ws.Range("E142").Value = 74

$ws.Range("B143").Value = 443
$ws.Range("C143").Value = 1
$ws.Range("D143").Value = 427
$ws.Range("E143").Value = 9
$ws.Range("H143").Value = 7

$ws.Range("D144").Value = 211
$ws.Range("E144").Value = 218
$ws.Range("H144").Value = 13

$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

$ws.Range("D202").Value = 16
$ws.Range("H202").Value = 2

$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
